$d = $word.ActiveDocument

function Set-ParagraphXml($para, [string]$innerXml) {
    $pkg = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$innerXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    [void]$para.Range.InsertXML($pkg)
}

# Locate the paragraph that currently reads exactly
# "They should test all the paths and possible outcomes"
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "They should test all the paths and possible outcomes") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not find the 'They should test...' paragraph"
}

$newInadequateParaXml = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>The test cases for these 3 classes are clearly inadequate since for all 3 of these classes 0% of the cases were covered. Instead, t</w:t></w:r>' +
    '<w:r><w:t>hey should test all the paths and possible outcomes</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> to try to reach maximum coverage of the classes. </w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

Set-ParagraphXml $target $newInadequateParaXml

# Locate the paragraph that currently reads exactly "83.5%" (this is the
# paragraph that used to hold the _GoBack bookmark; it must lose it now
# that the bookmark moved to the paragraph above).
$target2 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "83.5%") {
        $target2 = $p
        break
    }
}

if ($null -eq $target2) {
    throw "Could not find the '83.5%' paragraph"
}

$newScoreParaXml = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>83.5</w:t></w:r>' +
    '<w:r><w:t>%</w:t></w:r>' +
    '</w:p>'

Set-ParagraphXml $target2 $newScoreParaXml

Write-Output "Edit applied."
